# Adição da coluna quantidade a encomendar
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the TME row (row 3) - "OBTBT0001C" / "TME" duplicate offer line.
$ws.Rows.Item(3).Delete()

# 2) Remove the trailing "OCEBB0015S" row which is no longer part of the list.
#    After the previous delete, that row is now row 10.
$ws.Rows.Item(10).Delete()

# 3) Update the remaining offer rows (now rows 3-9) with the new best-quote data
#    for columns E (Fornecedor), F (Preço), G (MOQ), H (Prazo (dias)), I (Total_preço).
$rowData = @(
    @{ Row = 3; E = "Ariat";   F = 0.116;  G = 500;  H = 3; I = 58 },
    @{ Row = 4; E = "Ariat";   F = 0.139;  G = 500;  H = 3; I = 69.5 },
    @{ Row = 5; E = "Ariat";   F = 1.7;    G = 500;  H = 3; I = 850 },
    @{ Row = 6; E = "Ariat";   F = 0.72;   G = 1000; H = 3; I = 720 },
    @{ Row = 7; E = "Ariat";   F = 0.24;   G = 1000; H = 3; I = 240 },
    @{ Row = 8; E = "Ariat";   F = 0.023;  G = 500;  H = 3; I = 11.5 },
    @{ Row = 9; E = "Simento"; F = 0.0098; G = 500;  H = 3; I = 39.2 }
)

foreach ($r in $rowData) {
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
}

# 4) Add the new column J ("Qt a encomendar") with header + values for rows 1-9.
$ws.Cells.Item(1, 10).Value = "Qt a encomendar"
# Copy the header formatting (bold, centered, bordered) from I1 onto the new J1 header.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

$jValues = @{ 2 = 0; 3 = 500; 4 = 500; 5 = 500; 6 = 1000; 7 = 1000; 8 = 500; 9 = 4000 }
foreach ($rowNum in $jValues.Keys) {
    $ws.Cells.Item($rowNum, 10).Value = $jValues[$rowNum]
}
